# Adapt tests to control version
# Adds a "version" column (C) to the "settings" sheet, with value 1,
# and makes the settings sheet the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# Add new "version" column header and value
$ws.Range("C1").Value = "version"
$ws.Range("C2").Value = 1

# Select cell C3 on the settings sheet (mirrors prior A3 selection pattern)
$ws.Range("C3").Select()

# Make "settings" the active sheet/tab
$ws.Activate()
